$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.374260666666666
$ws.Range("H2").Value = 7.122781999999999
$ws.Range("I2").Value = 0.0276017086472712
$ws.Range("J2").Value = 0.0276017086472712
$ws.Range("M2").Value = 1.302860333333333
$ws.Range("N2").Value = 3.908581
$ws.Range("O2").Value = 0.9669439908960468
$ws.Range("P2").Value = 0.9669439908960467
$ws.Range("Q2").Value = 3.093330043593555
$ws.Range("R2").Value = 27.839970392342
$ws.Range("S2").Value = 0.02668930631494234
$ws.Range("T2").Value = 0.02668930631494234

# Row 3
$ws.Range("G3").Value = 2.374260666666666
$ws.Range("H3").Value = 7.122781999999999
$ws.Range("I3").Value = 0.0276017086472712
$ws.Range("J3").Value = 0.0276017086472712
$ws.Range("O3").Value = 0.008324674682103805
$ws.Range("P3").Value = 0.008324674682103805
$ws.Range("Q3").Value = 0.02663129047777777
$ws.Range("R3").Value = 0.2396816143
$ws.Range("S3").Value = 0.0002297752451587442
$ws.Range("T3").Value = 0.0002297752451587442

# Row 4
$ws.Range("G4").Value = 2.374260666666666
$ws.Range("H4").Value = 7.122781999999999
$ws.Range("I4").Value = 0.0276017086472712
$ws.Range("J4").Value = 0.0276017086472712
$ws.Range("M4").Value = 0.033323
$ws.Range("N4").Value = 0.099969
$ws.Range("O4").Value = 0.02473133442184949
$ws.Range("P4").Value = 0.02473133442184949
$ws.Range("Q4").Value = 0.07911748819533332
$ws.Range("R4").Value = 0.7120573937579999
$ws.Range("S4").Value = 0.0006826270871701189
$ws.Range("T4").Value = 0.0006826270871701189

# Row 5
$ws.Range("I5").Value = 0.8942818522422411
$ws.Range("J5").Value = 0.8942818522422411
$ws.Range("M5").Value = 1.302860333333333
$ws.Range("N5").Value = 3.908581
$ws.Range("O5").Value = 0.9669439908960468
$ws.Range("P5").Value = 0.9669439908960467
$ws.Range("Q5").Value = 100.2223795755812
$ws.Range("R5").Value = 902.001416180231
$ws.Range("S5").Value = 0.8647204631930214
$ws.Range("T5").Value = 0.8647204631930213

# Row 6
$ws.Range("I6").Value = 0.8942818522422411
$ws.Range("J6").Value = 0.8942818522422411
$ws.Range("O6").Value = 0.008324674682103805
$ws.Range("P6").Value = 0.008324674682103805
$ws.Range("S6").Value = 0.00744460549402588
$ws.Range("T6").Value = 0.00744460549402588

# Row 7
$ws.Range("I7").Value = 0.8942818522422411
$ws.Range("J7").Value = 0.8942818522422411
$ws.Range("M7").Value = 0.033323
$ws.Range("N7").Value = 0.099969
$ws.Range("O7").Value = 0.02473133442184949
$ws.Range("P7").Value = 0.02473133442184949
$ws.Range("Q7").Value = 2.563367898424334
$ws.Range("R7").Value = 23.070311085819
$ws.Range("S7").Value = 0.02211678355519385
$ws.Range("T7").Value = 0.02211678355519385

# Row 8
$ws.Range("G8").Value = 6.625048
$ws.Range("H8").Value = 19.875144
$ws.Range("I8").Value = 0.07701877356495823
$ws.Range("J8").Value = 0.07701877356495825
$ws.Range("M8").Value = 1.302860333333333
$ws.Range("N8").Value = 3.908581
$ws.Range("O8").Value = 0.9669439908960468
$ws.Range("P8").Value = 0.9669439908960467
$ws.Range("Q8").Value = 8.631512245629333
$ws.Range("R8").Value = 77.68361021066399
$ws.Range("S8").Value = 0.07447284028481967
$ws.Range("T8").Value = 0.07447284028481967

# Row 9
$ws.Range("G9").Value = 6.625048
$ws.Range("H9").Value = 19.875144
$ws.Range("I9").Value = 0.07701877356495823
$ws.Range("J9").Value = 0.07701877356495825
$ws.Range("O9").Value = 0.008324674682103805
$ws.Range("P9").Value = 0.008324674682103805
$ws.Range("Q9").Value = 0.07431095506666666
$ws.Range("R9").Value = 0.6687985956
$ws.Range("S9").Value = 0.0006411562343428935
$ws.Range("T9").Value = 0.0006411562343428937

# Row 10
$ws.Range("G10").Value = 6.625048
$ws.Range("H10").Value = 19.875144
$ws.Range("I10").Value = 0.07701877356495823
$ws.Range("J10").Value = 0.07701877356495825
$ws.Range("M10").Value = 0.033323
$ws.Range("N10").Value = 0.099969
$ws.Range("O10").Value = 0.02473133442184949
$ws.Range("P10").Value = 0.02473133442184949
$ws.Range("Q10").Value = 0.220766474504
$ws.Range("R10").Value = 1.986898270536
$ws.Range("S10").Value = 0.001904777045795683
$ws.Range("T10").Value = 0.001904777045795683

# Row 11
$ws.Range("G11").Value = 0.09441966666666667
$ws.Range("H11").Value = 0.283259
$ws.Range("I11").Value = 0.001097665545529457
$ws.Range("J11").Value = 0.001097665545529457
$ws.Range("M11").Value = 1.302860333333333
$ws.Range("N11").Value = 3.908581
$ws.Range("O11").Value = 0.9669439908960468
$ws.Range("P11").Value = 0.9669439908960467
$ws.Range("Q11").Value = 0.1230156383865555
$ws.Range("R11").Value = 1.107140745479
$ws.Range("S11").Value = 0.001061381103263339
$ws.Range("T11").Value = 0.001061381103263339

# Row 12
$ws.Range("G12").Value = 0.09441966666666667
$ws.Range("H12").Value = 0.283259
$ws.Range("I12").Value = 0.001097665545529457
$ws.Range("J12").Value = 0.001097665545529457
$ws.Range("O12").Value = 0.008324674682103805
$ws.Range("P12").Value = 0.008324674682103805
$ws.Range("Q12").Value = 0.001059073927777778
$ws.Range("R12").Value = 0.009531665349999999
$ws.Range("S12").Value = [double]"9.137708576286727E-06"
$ws.Range("T12").Value = [double]"9.137708576286727E-06"

# Row 13
$ws.Range("G13").Value = 0.09441966666666667
$ws.Range("H13").Value = 0.283259
$ws.Range("I13").Value = 0.001097665545529457
$ws.Range("J13").Value = 0.001097665545529457
$ws.Range("M13").Value = 0.033323
$ws.Range("N13").Value = 0.099969
$ws.Range("O13").Value = 0.02473133442184949
$ws.Range("P13").Value = 0.02473133442184949
$ws.Range("Q13").Value = 0.003146346552333333
$ws.Range("R13").Value = 0.028317118971
$ws.Range("S13").Value = [double]"2.714673368983084E-05"
$ws.Range("T13").Value = [double]"2.714673368983084E-05"
